# Append the new EUR->ARS quote row (row 26) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these as literal text (not auto-parsed dates/times/numbers),
# then reset the style so no extra formatting/style index lingers on the cells.
$ws.Cells.Item(26, 1).NumberFormat = "@"
$ws.Cells.Item(26, 1).Value = "2025-09-18"
$ws.Cells.Item(26, 1).Style = "Normal"

$ws.Cells.Item(26, 2).NumberFormat = "@"
$ws.Cells.Item(26, 2).Value = "21:23:24"
$ws.Cells.Item(26, 2).Style = "Normal"

$ws.Cells.Item(26, 3).NumberFormat = "@"
$ws.Cells.Item(26, 3).Value = "1.00 EUR = 1,760.9921"
$ws.Cells.Item(26, 3).Style = "Normal"
